# Chile coordinates + fix of csv delimiter
# Fill in the "Region" subdivision category (column E) for Chile's rows,
# which had been dropped because the source CSV used ";" as delimiter,
# and populate the Latitude (G) / Longitude (H) coordinates that were
# missing for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Latitude, Longitude)
$coords = @{
    88  = @(-46.035339,            -73.036434999999997)
    89  = @(-23.446308999999999,   -68.998489000000006)
    90  = @(-18.539467999999999,   -69.714607999999998)
    91  = @(-38.782228000000003,   -72.543251999999995)
    92  = @(-27.648325,            -70.432531999999995)
    93  = @(-37.225611999999998,   -73.108534000000006)
    94  = @(-30.540376999999999,   -70.967738999999995)
    95  = @(-34.434587999999998,   -71.154459000000003)
    96  = @(-42.046638000000002,   -73.008763999999999)
    97  = @(-39.922539,            -72.588984999999994)
    98  = @(-51.902416000000002,   -73.244017999999997)
    99  = @(-35.500971999999997,   -71.727129000000005)
    100 = @(-36.602809999999998,   -72.073119000000005)
    101 = @(-33.478729000000001,   -70.590025999999995)
    102 = @(-20.099081000000002,   -69.456920999999994)
    103 = @(-32.740869000000004,   -71.404539)
}

foreach ($row in 88..103) {
    $ws.Range("E$row").Value = "Region"

    $lat = $coords[$row][0]
    $lon = $coords[$row][1]
    $ws.Range("G$row").Value = $lat
    $ws.Range("H$row").Value = $lon
}

# Restore the current view position/selection as captured in the saved file
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("F101").Select()
